$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activities")

# New header for column J (copy header formatting from an existing header cell)
$ws.Range("J1").Value = "Average Heart rate (b/min)"
$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Full replacement data set, rows 2-12, columns A-J
$data = @(
    @("Evening Run",   15.5,  "01:34:05", 5645, 5850, "Run", "2025-06-29T18:01:37Z", "06:04", 14.348911104, 149.3),
    @("Intervallit",   10.27, "00:46:54", 2814, 2867, "Run", "2025-06-26T17:15:08Z", "04:34", 10.87916544,  168.5),
    @("Säbä",          7.42,  "01:35:30", 5730, 6304, "Run", "2025-06-23T19:03:11Z", "12:52", 11.39415552,  144.2),
    @("Afternoon Run", 7.84,  "00:36:20", 2180, 2187, "Run", "2025-06-22T15:48:10Z", "04:38", 7.7248512,    166.5),
    @("Afternoon Run", 13.63, "01:21:54", 4914, 4934, "Run", "2025-06-18T16:37:50Z", "06:01", 6.9201792,    152.5),
    @("Evening Run",   4.7,   "01:16:24", 4584, 4584, "Run", "2025-06-16T19:27:54Z", "16:15", 8.78701824,   142.4),
    @("Afternoon Run", 3.46,  "00:20:36", 1236, 1242, "Run", "2025-06-15T13:01:30Z", "05:57", 6.66268416,   $null),
    @("Säbä",          4.58,  "01:12:11", 4331, 4331, "Run", "2025-06-09T19:33:18Z", "15:46", 11.71602432,  $null),
    @("Afternoon Run", 7.19,  "00:47:31", 2851, 2909, "Run", "2025-06-04T16:29:24Z", "06:37", 6.46956288,   $null),
    @("Intervalli",    9.54,  "00:44:31", 2671, 2806, "Run", "2025-06-03T17:53:51Z", "04:40", 10.39636224,  $null),
    @("Lunch Run",     16.28, "01:44:35", 6275, 6876, "Run", "2025-06-01T12:37:04Z", "06:25", 7.892222976,  $null)
)

$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        $cell = $ws.Cells.Item($r, $c + 1)
        if ($null -eq $val) {
            $cell.Value = ""
        } else {
            $cell.Value = $val
        }
    }
    $r++
}
